$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Settings")
# B12 in the ORIGINAL sheet had style s="8" (quotePrefix text). Copy its format to B9.
$ws1.Range("B12").Copy()
$ws1.Range("B9").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("B9").Value = "16:30"
